# Add 2022-Q3 data:
#  - insert a new "2022-Q3" worksheet right after "总计", built from a copy of
#    the "2022-Q2" sheet (so it inherits the exact header/column-A styling),
#    trimmed down to its own 3 data rows and re-populated with 2022-Q3 figures.
#  - update the "总计" (summary) sheet with a new top row for 2022-Q3 and
#    shift the rest of the quarters down by one row.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet from a copy of "2022-Q2" (keeps styles).
# ---------------------------------------------------------------------------
$q2Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Drop rows 5-10 (the copy has 9 data rows, 2022-Q3 only needs 3).
$q3Sheet.Range("A5:H10").EntireRow.Delete()

# Row 2
$q3Sheet.Cells.Item(2, 2).Value = "009010"
$q3Sheet.Cells.Item(2, 3).Value = "华夏兴阳一年持有期混合"
$q3Sheet.Cells.Item(2, 4).Value = "26.58"
$q3Sheet.Cells.Item(2, 5).Value = "88.59"
$q3Sheet.Cells.Item(2, 6).Value = "2.90"
$q3Sheet.Cells.Item(2, 7).Value = "0.7708"
$q3Sheet.Cells.Item(2, 8).Value = 8

# Row 3
$q3Sheet.Cells.Item(3, 2).Value = "008704"
$q3Sheet.Cells.Item(3, 3).Value = "广发高股息优享混合A"
$q3Sheet.Cells.Item(3, 4).Value = "2.33"
$q3Sheet.Cells.Item(3, 5).Value = "92.85"
$q3Sheet.Cells.Item(3, 6).Value = "5.68"
$q3Sheet.Cells.Item(3, 7).Value = "0.1323"
$q3Sheet.Cells.Item(3, 8).Value = 4

# Row 4
$q3Sheet.Cells.Item(4, 2).Value = "008705"
$q3Sheet.Cells.Item(4, 3).Value = "广发高股息优享混合C"
$q3Sheet.Cells.Item(4, 4).Value = "0.69"
$q3Sheet.Cells.Item(4, 5).Value = "92.85"
$q3Sheet.Cells.Item(4, 6).Value = "5.68"
$q3Sheet.Cells.Item(4, 7).Value = "0.0392"
$q3Sheet.Cells.Item(4, 8).Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: new top row + shift the rest down one.
# ---------------------------------------------------------------------------

# Clone the A-column style (bold/centered/bordered) onto the new last row (8)
# before writing values, so it keeps the same formatting as every other
# row's "A" cell.
$totalSheet.Range("A2").Copy($totalSheet.Range("A8"))

$summary = @(
    @("2022-Q3", 3,  "0.9399999999999999"),
    @("2022-Q2", 9,  "2.9"),
    @("2022-Q1", 7,  "2.3"),
    @("2021-Q4", 7,  "0.84"),
    @("2021-Q3", 22, "7.09"),
    @("2021-Q2", 23, "3.82"),
    @("2021-Q1", 13, "0.9399999999999999")
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $r = $i + 2
    $entry = $summary[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $entry[0]
    $totalSheet.Cells.Item($r, 3).Value = $entry[1]
    $totalSheet.Cells.Item($r, 4).Value = [double]$entry[2]
}
